# corrección de error behavior en los test case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two "NroSiniestro" test-case values (F2, F3) with corrected numbers,
# preserving the original trailing whitespace pattern. A leading apostrophe forces
# Excel to keep the value as text (preserving the leading zero / trailing spaces)
# instead of converting it to a number, matching how the cells were already formatted.
$ws.Range("F3").Value = "'0420172010222  "
$ws.Range("F2").Value = "'0420194406906 "

# Update the active selection from F4 to G4
$ws.Range("G4").Select()
